# Added functionality to add name and number at the top
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Name" / "Number" columns with header labels
$ws.Range("F1").Value = "Name"
$ws.Range("G1").Value = "Number"

# First item's name and number
$ws.Range("F2").Value = "Nishika"
$ws.Range("G2").Value = 9871793958

# Match the font used for the new header/data cells
$ws.Range("F1:G2").Font.ThemeColor = 1

# Tidy up packing descriptions (KGS BAG -> KG/BAG)
$ws.Range("E2").Value = "50 KG/BAG"
$ws.Range("E3").Value = "50 KG/BAG"
$ws.Range("E4").Value = "25 KG/BAG"
